# Generate Report for Handoff
# Update the "Latest Handoff Datetime" column (D) for the b677b96d file row (row 3)
# on both language sheets, reflecting a newly-generated handoff report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-02-18 06:26:22"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-02-18 06:26:35"
